$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used as a staging area so that values which look numeric
# (e.g. "303.32", "4.89%") are written back as literal TEXT, matching the
# original t="inlineStr" cells, without leaving a NumberFormat/style trace
# on the destination cells (PasteSpecial values-only only copies value+type).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$helper.Value = "303.32"
$helper.Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$helper.Value = "4.89%"
$helper.Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4163) | Out-Null
$helper.Value = "34.95"
$helper.Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$helper.Value = "12.55%"
$helper.Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4163) | Out-Null
$helper.Value = "5.165"
$helper.Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$helper.Value = "4.30%"
$helper.Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.07819"
$helper.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$helper.Value = "6.25%"
$helper.Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4163) | Out-Null
$helper.Value = "2.286"
$helper.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$helper.Value = "-2.61%"
$helper.Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4163) | Out-Null
$helper.Value = "8.056"
$helper.Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$helper.Value = "4.17%"
$helper.Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4163) | Out-Null
$helper.Value = "3.988"
$helper.Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$helper.Value = "7.02%"
$helper.Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.9228"
$helper.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$helper.Value = "1.21%"
$helper.Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.1007"
$helper.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$helper.Value = "8.66%"
$helper.Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.1839"
$helper.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$helper.Value = "7.90%"
$helper.Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.08528"
$helper.Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$helper.Value = "4.40%"
$helper.Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.03375"
$helper.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$helper.Value = "8.40%"
$helper.Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.09911"
$helper.Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$helper.Value = "-0.63%"
$helper.Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.001481"
$helper.Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$helper.Value = "-1.41%"
$helper.Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4163) | Out-Null
$helper.Value = "2.79%"
$helper.Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.005746"
$helper.Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.47%"
$helper.Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.33%"
$helper.Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4163) | Out-Null
$helper.Value = "2.126"
$helper.Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.95%"
$helper.Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.3416"
$helper.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$helper.Value = "2.74%"
$helper.Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4163) | Out-Null
$helper.Value = "4.568"
$helper.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$helper.Value = "9.32%"
$helper.Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.2395"
$helper.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$helper.Value = "14.05%"
$helper.Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.001224"
$helper.Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.99%"
$helper.Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.004329"
$helper.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$helper.Value = "3.58%"
$helper.Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.0001302"
$helper.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.19%"
$helper.Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.0003400"
$helper.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.22%"
$helper.Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.01742"
$helper.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$helper.Value = "10.52%"
$helper.Copy() | Out-Null
$ws.Range("E39").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.04740"
$helper.Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$helper.Value = "6.08%"
$helper.Copy() | Out-Null
$ws.Range("E40").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.007684"
$helper.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.1414"
$helper.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$helper.Value = "6.29%"
$helper.Copy() | Out-Null
$ws.Range("E42").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.007373"
$helper.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$helper.Value = "-25.79%"
$helper.Copy() | Out-Null
$ws.Range("E43").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.002212"
$helper.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$helper.Value = "-1.24%"
$helper.Copy() | Out-Null
$ws.Range("E44").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.01006"
$helper.Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$helper.Value = "14.57%"
$helper.Copy() | Out-Null
$ws.Range("E45").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.00006065"
$helper.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$helper.Value = "-0.69%"
$helper.Copy() | Out-Null
$ws.Range("E46").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.00000000750"
$helper.Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.05%"
$helper.Copy() | Out-Null
$ws.Range("E47").PasteSpecial(-4163) | Out-Null
$helper.Value = "3.879"
$helper.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$helper.Value = "58.64%"
$helper.Copy() | Out-Null
$ws.Range("E48").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.002690"
$helper.Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$helper.Value = "34.57%"
$helper.Copy() | Out-Null
$ws.Range("E49").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.00002100"
$helper.Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.05%"
$helper.Copy() | Out-Null
$ws.Range("E50").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.0002000"
$helper.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$helper.Value = "0.05%"
$helper.Copy() | Out-Null
$ws.Range("E51").PasteSpecial(-4163) | Out-Null

$helper.Clear() | Out-Null
$excel.CutCopyMode = $false

